$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (shifts existing rows 5-14 down to 6-15),
# copying the formatting from the row above (row 4) onto the new row.
$ws.Rows.Item(5).Insert(-4121, -4163)

# Populate the new row 5 with the Currency fields ("Měna dokladu" / "Currency" / "Currency")
$ws.Cells.Item(5, 1).Value = "Měna dokladu"
$ws.Cells.Item(5, 2).Value = "Currency"
$ws.Cells.Item(5, 3).Value = "Currency"

# Match the final on-screen selection recorded in the saved workbook
$ws.Range("C18").Select()
